$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume 1h (E) figures from the latest crypto data refresh.
# D-column values are forced to Text (matching the original inline-string cells) so that
# numeric-looking prices (e.g. "324.25") are not auto-converted to numbers by Excel,
# then formatting is cleared so the cell style stays identical to the untouched cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.184.37"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.99%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.532.59"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.43%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.25"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.95"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.528"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.59%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.559"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.70"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.50"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +11.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0829"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.88%  "

$ws.Range("E13").Value = "  +1.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.31"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.929.30"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.532.80"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.861"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "48.052.13"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.85%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.28"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.66%  "

$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0952"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.69"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.35"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.16"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +9.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.59"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.26"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.88%  "

$ws.Range("E27").Value = "  -0.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.18"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("E29").Value = "  +3.97%  "

$ws.Range("E30").Value = "  -3.71%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.77"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.80"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.93"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.42"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("E35").Value = "  -0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0795"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.69%  "

$ws.Range("E37").Value = "  +0.90%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.77"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.37%  "

$ws.Range("E39").Value = "  +1.05%  "

$ws.Range("E40").Value = "  +0.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.35"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.30"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.22%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.19"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.61%  "

$ws.Range("E44").Value = "  +0.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.014.75"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.00%  "

$ws.Range("E46").Value = "  +3.34%  "

$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("E48").Value = "  +5.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.18"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.65%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.27"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.84"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.37%  "
